# Compare-SmartData: use consistent variable naming ("criteria" instead of the
# misspelled "criterai") and make the output clearer by explicitly flagging
# every row that has no pass/fail criteria as a WARNING instead of leaving the
# cell blank.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that previously had a blank "criteria" column now get an explicit
# WARNING marker so the reader knows those rows were intentionally skipped
# rather than forgotten.
$warningRows = @(6, 7, 8, 9, 10, 11, 14, 15, 16, 19, 20, 22, 23, 26)
foreach ($r in $warningRows) {
    $ws.Cells.Item($r, 4).Value = "WARNING"
}

# Fix the misspelled column header.
$ws.Range("D1").Value = "criteria"

# Reflect the new content widths (column A grew, and column D now stands on
# its own instead of being merged with column E's width).
$ws.Columns.Item(1).ColumnWidth = 10.166666666666666
$ws.Columns.Item(4).ColumnWidth = 8.8

# Move the active selection to the header cell that was just corrected.
[void]$ws.Range("D1").Select()
